# Update the "historias de usuario" (user stories) worksheet:
# - fix spelling/accent issues in several existing story fragments
#   (e.g. "pagina" -> "página", "cuantos" -> "cuántos", "boton" -> "botón")
# - replace several story fragments with new/updated wording (new cat-photo
#   related stories, reworded "home page"/"search bar"/"login tab" stories, etc.)
# - row 21: the "tipo de usuario" cell is corrected back to " usuario "

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 'subir fotos '
$ws.Range("E3").Value = 'los demás las puedan ver'

$ws.Range("D4").Value = 'ver fotos '
$ws.Range("E4").Value = 'pueda apreciar la ternura de los animales'

$ws.Range("D5").Value = 'tener una cuenta '
$ws.Range("E5").Value = 'tenga una identidad'

$ws.Range("D6").Value = 'hacer login '
$ws.Range("E6").Value = 'todo lo que haga se vincule a mi cuenta '

$ws.Range("D7").Value = 'dar like a fotos '
$ws.Range("E7").Value = 'los demás puedan ver que le di like '

$ws.Range("D8").Value = 'ver cuántos likes tiene una foto '
$ws.Range("E8").Value = 'pueda saber qué fotos son populares '

$ws.Range("D9").Value = 'seguir a otros usuarios '
$ws.Range("E9").Value = 'pueda ver sus fotos en mi página de inicio '

$ws.Range("D10").Value = 'ser seguido por otros usuarios '
$ws.Range("E10").Value = 'puedan ver mis fotos en su página de inicio '

$ws.Range("D11").Value = 'tener una página de inicio '
$ws.Range("E11").Value = 'pueda ver las fotos de los usuarios que sigo '

$ws.Range("D12").Value = 'poder hacer log out '
$ws.Range("E12").Value = 'mi cuenta se cierre '

$ws.Range("D13").Value = 'tener una pestaña de log in '
$ws.Range("E13").Value = 'me dirija a mi página de inicio '

$ws.Range("D14").Value = 'tener una lista de las personas que sigo '
$ws.Range("E14").Value = 'pueda cambiar esa lista y se actualicen los datos '

$ws.Range("D15").Value = 'tener un perfil '
$ws.Range("E15").Value = 'puedan ver mis fotos, seguirme, ver a quién sigo y quién me sigue '

$ws.Range("D16").Value = 'dar click en una foto y poder verla más grande'
$ws.Range("E16").Value = 'pueda ver mejor la foto '

$ws.Range("D17").Value = 'tener una barra de búsqueda '
$ws.Range("E17").Value = 'al dar "buscar " me aparezcan fotos o perfiles relacionados'

$ws.Range("D18").Value = 'saber cuántas publicaciones he realizado'
$ws.Range("E18").Value = 'la gente y yo podamos ver qué tan activo soy'

$ws.Range("D19").Value = 'poder acceder desde mi pc o celular '
$ws.Range("E19").Value = 'sea cómodo navegar en cualquiera de ellos '

$ws.Range("D20").Value = 'dar click en un botón '
$ws.Range("E20").Value = 'pueda acceder a mi perfil '

$ws.Range("C21").Value = ' usuario '
$ws.Range("D21").Value = 'dar click en un botón '
$ws.Range("E21").Value = 'pueda ir a mi página de inicio '

$ws.Range("D22").Value = 'ver una campanita '
$ws.Range("E22").Value = ' me indique si tengo notificaciones'

$ws.Range("D23").Value = 'ver cuántos seguidores tengo '
$ws.Range("E23").Value = 'esta cantidad aparezca en mi perfil '

# Restore selection to the cell last active when the file was saved.
$ws.Range("D23").Select()
